$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to be treated as plain text so Excel does not
    # auto-convert numeric-looking / percent-looking strings into
    # numbers with a different number format (which would alter
    # styles and cell value types). After writing the value we put
    # the style back to "Normal" so no residual formatting remains.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2,4) "321.93"
Set-TextValue $ws.Cells.Item(2,5) "6.70%"
Set-TextValue $ws.Cells.Item(2,7) "17"
Set-TextValue $ws.Cells.Item(3,4) "49.22"
Set-TextValue $ws.Cells.Item(3,5) "11.47%"
Set-TextValue $ws.Cells.Item(3,7) "17"
Set-TextValue $ws.Cells.Item(4,4) "5.343"
Set-TextValue $ws.Cells.Item(4,5) "5.08%"
Set-TextValue $ws.Cells.Item(4,7) "17"
Set-TextValue $ws.Cells.Item(5,4) "0.08048"
Set-TextValue $ws.Cells.Item(5,5) "4.57%"
Set-TextValue $ws.Cells.Item(5,7) "17"
Set-TextValue $ws.Cells.Item(6,4) "4.604"
Set-TextValue $ws.Cells.Item(6,5) "4.13%"
Set-TextValue $ws.Cells.Item(6,7) "17"
Set-TextValue $ws.Cells.Item(7,4) "1.396"
Set-TextValue $ws.Cells.Item(7,5) "34.45%"
Set-TextValue $ws.Cells.Item(7,7) "17"
Set-TextValue $ws.Cells.Item(8,4) "1.641"
Set-TextValue $ws.Cells.Item(8,5) "1.57%"
Set-TextValue $ws.Cells.Item(8,7) "17"
Set-TextValue $ws.Cells.Item(9,4) "0.1284"
Set-TextValue $ws.Cells.Item(9,5) "0.77%"
Set-TextValue $ws.Cells.Item(9,7) "17"
Set-TextValue $ws.Cells.Item(10,4) "0.1965"
Set-TextValue $ws.Cells.Item(10,5) "5.32%"
Set-TextValue $ws.Cells.Item(10,7) "17"
Set-TextValue $ws.Cells.Item(11,4) "0.09576"
Set-TextValue $ws.Cells.Item(11,5) "3.61%"
Set-TextValue $ws.Cells.Item(11,7) "17"
Set-TextValue $ws.Cells.Item(12,4) "0.04704"
Set-TextValue $ws.Cells.Item(12,5) "12.74%"
Set-TextValue $ws.Cells.Item(12,7) "17"
Set-TextValue $ws.Cells.Item(13,4) "0.1044"
Set-TextValue $ws.Cells.Item(13,5) "-0.51%"
Set-TextValue $ws.Cells.Item(13,7) "17"
Set-TextValue $ws.Cells.Item(14,5) "3.02%"
Set-TextValue $ws.Cells.Item(14,7) "17"
Set-TextValue $ws.Cells.Item(15,4) "0.04198"
Set-TextValue $ws.Cells.Item(15,5) "-0.12%"
Set-TextValue $ws.Cells.Item(15,7) "17"
Set-TextValue $ws.Cells.Item(16,4) "0.005792"
Set-TextValue $ws.Cells.Item(16,5) "0.24%"
Set-TextValue $ws.Cells.Item(16,7) "17"
Set-TextValue $ws.Cells.Item(17,4) "3.341"
Set-TextValue $ws.Cells.Item(17,5) "-0.15%"
Set-TextValue $ws.Cells.Item(17,7) "17"
Set-TextValue $ws.Cells.Item(18,4) "2.448"
Set-TextValue $ws.Cells.Item(18,5) "5.02%"
Set-TextValue $ws.Cells.Item(18,7) "17"
Set-TextValue $ws.Cells.Item(19,4) "0.3510"
Set-TextValue $ws.Cells.Item(19,5) "4.69%"
Set-TextValue $ws.Cells.Item(19,7) "17"
Set-TextValue $ws.Cells.Item(20,4) "7.999"
Set-TextValue $ws.Cells.Item(20,5) "-7.79%"
Set-TextValue $ws.Cells.Item(20,7) "17"
Set-TextValue $ws.Cells.Item(21,4) "0.1373"
Set-TextValue $ws.Cells.Item(21,5) "-1.92%"
Set-TextValue $ws.Cells.Item(21,7) "17"
Set-TextValue $ws.Cells.Item(22,4) "0.3084"
Set-TextValue $ws.Cells.Item(22,5) "-2.96%"
Set-TextValue $ws.Cells.Item(22,7) "17"
Set-TextValue $ws.Cells.Item(23,4) "0.001311"
Set-TextValue $ws.Cells.Item(23,5) "1.91%"
Set-TextValue $ws.Cells.Item(23,7) "17"
Set-TextValue $ws.Cells.Item(24,4) "0.004267"
Set-TextValue $ws.Cells.Item(24,5) "-4.57%"
Set-TextValue $ws.Cells.Item(24,7) "17"
Set-TextValue $ws.Cells.Item(25,4) "0.0001345"
Set-TextValue $ws.Cells.Item(25,5) "-0.35%"
Set-TextValue $ws.Cells.Item(25,7) "17"
Set-TextValue $ws.Cells.Item(26,4) "0.0003527"
Set-TextValue $ws.Cells.Item(26,7) "17"
Set-TextValue $ws.Cells.Item(27,7) "17"
Set-TextValue $ws.Cells.Item(28,7) "17"
Set-TextValue $ws.Cells.Item(29,7) "17"
Set-TextValue $ws.Cells.Item(30,7) "17"
Set-TextValue $ws.Cells.Item(31,7) "17"
Set-TextValue $ws.Cells.Item(32,7) "17"
Set-TextValue $ws.Cells.Item(33,7) "17"
Set-TextValue $ws.Cells.Item(34,7) "17"
Set-TextValue $ws.Cells.Item(35,7) "17"
Set-TextValue $ws.Cells.Item(36,7) "17"
Set-TextValue $ws.Cells.Item(37,7) "17"
Set-TextValue $ws.Cells.Item(38,4) "0.02742"
Set-TextValue $ws.Cells.Item(38,5) "9.71%"
Set-TextValue $ws.Cells.Item(38,7) "17"
Set-TextValue $ws.Cells.Item(39,4) "0.06200"
Set-TextValue $ws.Cells.Item(39,5) "17.39%"
Set-TextValue $ws.Cells.Item(39,7) "17"
Set-TextValue $ws.Cells.Item(40,4) "0.01083"
Set-TextValue $ws.Cells.Item(40,5) "82.12%"
Set-TextValue $ws.Cells.Item(40,7) "17"
Set-TextValue $ws.Cells.Item(41,4) "0.008044"
Set-TextValue $ws.Cells.Item(41,5) "4.10%"
Set-TextValue $ws.Cells.Item(41,7) "17"
Set-TextValue $ws.Cells.Item(42,4) "0.1465"
Set-TextValue $ws.Cells.Item(42,5) "8.52%"
Set-TextValue $ws.Cells.Item(42,7) "17"
Set-TextValue $ws.Cells.Item(43,4) "0.007877"
Set-TextValue $ws.Cells.Item(43,5) "7.08%"
Set-TextValue $ws.Cells.Item(43,7) "17"
Set-TextValue $ws.Cells.Item(44,4) "0.008639"
Set-TextValue $ws.Cells.Item(44,5) "14.14%"
Set-TextValue $ws.Cells.Item(44,7) "17"
Set-TextValue $ws.Cells.Item(45,4) "0.3509"
Set-TextValue $ws.Cells.Item(45,5) "16.99%"
Set-TextValue $ws.Cells.Item(45,7) "17"
Set-TextValue $ws.Cells.Item(46,4) "0.00006746"
Set-TextValue $ws.Cells.Item(46,5) "1.42%"
Set-TextValue $ws.Cells.Item(46,7) "17"
Set-TextValue $ws.Cells.Item(47,4) "0.00000000747"
Set-TextValue $ws.Cells.Item(47,5) "-0.39%"
Set-TextValue $ws.Cells.Item(47,7) "17"
Set-TextValue $ws.Cells.Item(48,4) "0.05518"
Set-TextValue $ws.Cells.Item(48,5) "23.60%"
Set-TextValue $ws.Cells.Item(48,7) "17"
Set-TextValue $ws.Cells.Item(49,4) "0.003986"
Set-TextValue $ws.Cells.Item(49,5) "-5.10%"
Set-TextValue $ws.Cells.Item(49,7) "17"
Set-TextValue $ws.Cells.Item(50,4) "0.00002093"
Set-TextValue $ws.Cells.Item(50,5) "-0.39%"
Set-TextValue $ws.Cells.Item(50,7) "17"
Set-TextValue $ws.Cells.Item(51,4) "0.0001993"
Set-TextValue $ws.Cells.Item(51,5) "-0.39%"
Set-TextValue $ws.Cells.Item(51,7) "17"
